# Wrapped up week1 changes
# - Rewrote the "Merge Intervals" problem description / solution text
# - Rewrote the "Group Anagrams" problem description text (added example)
# - Adjusted row 10 height to fit the new (longer) content
# - Updated the active selection / scroll position to reflect where the
#   author ended up after editing (near row 15 / C15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Problem Description for "Merge Intervals" (row 10, col C) ---
$c10 = 'Given an array of Intervals, merge any overlapping intervals and return the output array. ' + "`n" + `
       'Ex. Input = {{1,3}, {2, 6}, {8, 10}, {15, 18}} Output = {{1, 6}, {8, 10}, {15, 18}}'
$ws.Range("C10").Value2 = $c10

# --- New Solution text for "Merge Intervals" (row 10, col D) ---
$d10 = 'Sort the array using -> ex Arrays.sort(intervals, (arr1), (arr2) -> Interger.compare(arr1[0], arr[1])) ' + "`n" + `
       "Need to create a list of int[] because we don't know the size at runtime bc of merges.  Have a pointer to current_interval OUTSIDE of loop.  In loop, If you need to add an array to output, first set current_interval given interval array and add current_int to output. Merge:  current_int[1] = Math.max(currHigh, nextHigh) "
$ws.Range("D10").Value2 = $d10

# --- New Problem Description for "Group Anagrams" (row 11, col C) ---
$c11 = 'Given an array of Strings, return a list of arrays with where a given array is anagrams of eachother' + "`n" + `
       'Input -> strs = ["eat", "tea", "tan", "ate", "nat", "bat"]  ' + "`n" + `
       'Output -> [["bat"], ["nat", "tan"], ["ate", "eat", "tea"] '
$ws.Range("C11").Value2 = $c11

# Row 10 grew slightly taller to accommodate the new wording.
$ws.Rows.Item(10).RowHeight = 61.5

# Reflect the author's final cursor / scroll position when they saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C15").Select() | Out-Null
